# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: on sheet "Rules", cell B11 goes from the text "R40" to the
# text "1" (still a plain text value, not a number).
#
# A direct `Range.Value = "1"` gets auto-typed as a *number* by Excel's
# input parser (same as typing 1 into a General-formatted cell), which
# would store it with no `t="s"` shared-string typing. Forcing text with
# a leading apostrophe avoids that, but it also stamps a quote-prefix
# onto the cell's style, creating a brand-new style record and changing
# B11's `s=` (style) index - which should stay exactly as it was.
#
# Workaround: stage the text "1" in a scratch cell (forcing it to text
# via the apostrophe there), copy it, and paste *values only* onto B11.
# A values-only paste carries over the string-typed content without
# touching B11's existing style/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Value = "'1"

$scratch.Copy()
$target = $ws.Range("B11")
$target.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$scratch.Clear()
$excel.CutCopyMode = $false
